$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire second row (emailAddress "dpaul" / "password" + hyperlink),
# which shifts the "abcdef" row up to become row 2.
$ws.Rows.Item(2).Delete()

# The deleted row carried a mailto hyperlink; make sure none remain on the sheet.
$ws.Hyperlinks.Delete()

# Move the active selection (matches the target workbook state).
$ws.Range("E4").Select()
